$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 data updates (Creditos test fixture refresh):
#   usuario (D4), clave (E4), tipoPrestamo (O4), numeroPrestamo (P4), numeroCuenta (T4)
$ws.Range("D4").Value2 = "pruebauser01"
$ws.Range("E4").Value2 = "6789"
$ws.Range("O4").Value2 = "Prestamo personal"
$ws.Range("P4").Value2 = "29281023961"
$ws.Range("T4").Value2 = "406-182800-03"

# Writing new text above resets the cells' "quote-prefixed text" formatting
# (style shared with the rest of the row, e.g. C4/F4, s=5). Reapply it from
# an untouched neighbour so D4/E4 keep the original look.
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
